$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the 3 runs "ves the " + "filter" + " transform."
# into a single run "ves the filter transform." (leaving the
# preceding "...that dri" run untouched).
# -----------------------------------------------------------------
$rngFilter = $d.Content
$foundFilter = $rngFilter.Find.Execute("filter", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if (-not $foundFilter) { throw "Could not find 'filter'" }
$filterStart = $rngFilter.Start
$filterEnd = $rngFilter.End

# 1a) Touch the boundary between "filter" and " transform." (edit
#     starting exactly at $filterEnd) so the two runs merge into one
#     (force a real content change, corrected again in step 1c).
$rA = $d.Range($filterEnd, $filterEnd + 11)
$rA.Text = " transformX."

# 1b) Touch the boundary between "ves the " and the just-merged
#     "filter transformX." run (edit starting exactly at
#     $filterStart) so all three runs become one.
$rB = $d.Range($filterStart, $filterStart + "filter transformX.".Length)
$rB.Text = "filter transformX."

# 1c) Fix the temporary "X" marker -- strictly inside the merged run
#     (not at its start boundary) so no further merge is triggered.
$rngC = $d.Content
$foundC = $rngC.Find.Execute("transformX.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if (-not $foundC) { throw "Could not find 'transformX.'" }
$rC = $d.Range($rngC.Start, $rngC.Start + 10)
$rC.Text = "transform"

# -----------------------------------------------------------------
# Change 2 & 5: relocate the "_GoBack" bookmark.
#  - remove it from its old spot (inside "crucial", Design Assurance
#    section)
#  - remove the dangling AADL file/line reference text near
#    "SW.Impl" ("on line 70 in Producer_Filter_Consumer.aadl")
#  - re-add a (now empty) "_GoBack" bookmark right after "SW.Impl"
# -----------------------------------------------------------------
$bmGoBack = $d.Bookmarks("_GoBack")
$bmGoBack.Delete()

# Delete " on line 70 in Producer_Filter_Consumer.aadl" but leave the
# last letter of "SW.Impl" as part of the touched range so the extra
# duplicate proofErr spellEnd produced by Word's run-splitting
# collapses back down to a single one (matching real Word behaviour).
$rngDel = $d.Content
$foundDel = $rngDel.Find.Execute("l on line 70 in Producer_Filter_Consumer.aadl", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if (-not $foundDel) { throw "Could not find the aadl file reference text" }
$rngDel.Text = "l"

# Re-find "SW.Impl" (position may have shifted after the bookmark
# delete / text delete above) and drop the bookmark right after it.
$rngSW = $d.Content
$foundSW = $rngSW.Find.Execute("SW.Impl", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if (-not $foundSW) { throw "Could not find 'SW.Impl'" }
$swEnd = $rngSW.End
$insertPoint = $d.Range($swEnd, $swEnd)
$d.Bookmarks.Add("_GoBack", $insertPoint)

Write-Output "Done."
